$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.045.47"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "1.886.92"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.24"
$ws.Range("E5").Value = "  -3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -3.89%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.77"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07949"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9969"
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.84"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("D13").Value = "1.903.97"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.914"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.068"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.54"
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06551"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001022"
$ws.Range("E19").Value = "  -3.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "29.070.42"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.433"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("D26").Value = "2.125.99"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.28"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.082"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.499"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.30"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.034"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09314"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.517"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.287"
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06052"
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02221"
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.341"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.173"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5776"
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.10"
$ws.Range("E44").Value = "  -5.46%  "
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07517"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.290"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.01"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5452"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.902"
$ws.Range("E50").Value = "  -4.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.64"
$ws.Range("E51").Value = "  -2.28%  "
